# Applies cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.380.38"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.889.72"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.694"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0742"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "2.163.62"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.762"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "1.880.43"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "35.498.25"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.05%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "4.128.47"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +9.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0693"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0220"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("D45").Value = "1.294.33"
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.41%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0797"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.96%  "
$ws.Range("B48").Value = "Gas"
$ws.Range("C48").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.06%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("E51").Value = "  -6.05%  "
